$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right before the "总计" sheet,
#    matching the style/format of the existing "2021-Q4" sheet.
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Copy header-row formatting (bold + border) from 2021-Q4's header row.
$src.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

# Copy column-A (index) cell formatting down through rows 2-9.
$src.Range("A2").Copy()
$q1.Range("A2:A9").PasteSpecial(-4122)   # xlPasteFormats

# Columns B:G hold text-like values (fund codes / numbers-as-text), so
# force a Text number format before assigning the values so they are
# not auto-coerced into numeric cells.
$q1.Range("B2:G9").NumberFormat = "@"

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data rows
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "160325"
$q1.Range("C2").Value = "华夏创业板两年定期开放混合"
$q1.Range("D2").Value = "27.39"
$q1.Range("E2").Value = "90.77"
$q1.Range("F2").Value = "4.46"
$q1.Range("G2").Value = "1.2216"
$q1.Range("H2").Value = 3

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "001042"
$q1.Range("C3").Value = "华夏领先股票"
$q1.Range("D3").Value = "13.72"
$q1.Range("E3").Value = "93.46"
$q1.Range("F3").Value = "6.65"
$q1.Range("G3").Value = "0.9124"
$q1.Range("H3").Value = 3

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "506006"
$q1.Range("C4").Value = "汇添富科创板2年定期开放混合"
$q1.Range("D4").Value = "26.37"
$q1.Range("E4").Value = "91.69"
$q1.Range("F4").Value = "3.40"
$q1.Range("G4").Value = "0.8966"
$q1.Range("H4").Value = 10

$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "519029"
$q1.Range("C5").Value = "华夏稳增混合"
$q1.Range("D5").Value = "10.92"
$q1.Range("E5").Value = "92.99"
$q1.Range("F5").Value = "4.04"
$q1.Range("G5").Value = "0.4412"
$q1.Range("H5").Value = 8

$q1.Range("A6").Value = 4
$q1.Range("B6").Value = "014185"
$q1.Range("C6").Value = "招商专精特新股票A"
$q1.Range("D6").Value = "8.37"
$q1.Range("E6").Value = "30.94"
$q1.Range("F6").Value = "1.15"
$q1.Range("G6").Value = "0.0963"
$q1.Range("H6").Value = 10

$q1.Range("A7").Value = 5
$q1.Range("B7").Value = "014186"
$q1.Range("C7").Value = "招商专精特新股票C"
$q1.Range("D7").Value = "3.46"
$q1.Range("E7").Value = "30.94"
$q1.Range("F7").Value = "1.15"
$q1.Range("G7").Value = "0.0398"
$q1.Range("H7").Value = 10

$q1.Range("A8").Value = 6
$q1.Range("B8").Value = "710301"
$q1.Range("C8").Value = "富安达增强收益债券A"
$q1.Range("D8").Value = "0.61"
$q1.Range("E8").Value = "20.20"
$q1.Range("F8").Value = "0.96"
$q1.Range("G8").Value = "0.0059"
$q1.Range("H8").Value = 8

$q1.Range("A9").Value = 7
$q1.Range("B9").Value = "710302"
$q1.Range("C9").Value = "富安达增强收益债券C"
$q1.Range("D9").Value = "0.26"
$q1.Range("E9").Value = "20.20"
$q1.Range("F9").Value = "0.96"
$q1.Range("G9").Value = "0.0025"
$q1.Range("H9").Value = 8

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new top data row for
#    "2022-Q1" and push the existing rows down by one, preserving the
#    existing formatting on column A.
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# Create row 5 by copying the style of the existing row-4 index cell.
$tot.Range("A4").Copy()
$tot.Range("A5").PasteSpecial(-4122)   # xlPasteFormats

# Cascade the existing data down one row (bottom-up to avoid clobbering).
$tot.Range("B5").Value = $tot.Range("B4").Value()
$tot.Range("C5").Value = $tot.Range("C4").Value()
$tot.Range("D5").Value = $tot.Range("D4").Value()
$tot.Range("A5").Value = 3

$tot.Range("B4").Value = $tot.Range("B3").Value()
$tot.Range("C4").Value = $tot.Range("C3").Value()
$tot.Range("D4").Value = $tot.Range("D3").Value()
$tot.Range("A4").Value = 2

$tot.Range("B3").Value = $tot.Range("B2").Value()
$tot.Range("C3").Value = $tot.Range("C2").Value()
$tot.Range("D3").Value = $tot.Range("D2").Value()
$tot.Range("A3").Value = 1

$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 8
$tot.Range("D2").Value = 3.62
$tot.Range("A2").Value = 0
